$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5
# from serial date 45221 (2023-10-22) to 45224 (2023-10-25)
$ws.Range("C2").Value = 45224
$ws.Range("C3").Value = 45224
$ws.Range("C4").Value = 45224
$ws.Range("C5").Value = 45224
